$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Row 4: replace old text with new task text
$ws.Range("B4").Value = "Write Summary of Chap 1 (half)"

# Row 5: the day-counter formula moves up here (was in row 6), plus new task text
$ws.Range("A5").Formula = "=A3+1"
$ws.Range("B5").Value = "Complete Calculus 12.1 Assignments"

# Clear out the old A6/A7/A8 counter cells - that sequence now starts at row 9
$ws.Range("A6:A8").ClearContents()

# New task rows for Nov 2nd
$ws.Range("B6").Value = "Complete Calculus 12.2 Assignments"
$ws.Range("B7").Value = "Write Summary of Chap 1 (other half)"
$ws.Range("B8").Value = "Complete Calculus 14.1 Assignment"

# Counter sequence resumes at row 9, continuing down through row 26
$ws.Range("A9").Formula = "=A5+1"
$ws.Range("A10").Formula = "=A9+1"
$ws.Range("A25").Formula = "=A24+1"
$ws.Range("A26").Formula = "=A25+1"

$ws.Range("B9").Select()
